$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 268.57144
$ws.Range("I2").Value = 268.57144
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 268.57144
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -155.57144
# Row 39
$ws.Range("H39").Value = 317.8
$ws.Range("I39").Value = 69.2381
$ws.Range("J39").Value = 897.7778
$ws.Range("K39").Value = 207.7143
$ws.Range("L39").Value = 2693.3334
$ws.Range("M39").Value = 88.28569999999999
$ws.Range("N39").Value = -3285.3334
# Row 58
$ws.Range("H58").Value = 1253.4
$ws.Range("I58").Value = 87.625
$ws.Range("J58").Value = 2585.7144
$ws.Range("K58").Value = 262.875
$ws.Range("L58").Value = 7757.1432
$ws.Range("M58").Value = -112.875
$ws.Range("N58").Value = -8057.1432
# Row 69
$ws.Range("H69").Value = 3761
$ws.Range("I69").Value = 3013
$ws.Range("J69").Value = 4509
$ws.Range("K69").Value = 9039
$ws.Range("L69").Value = 13527
$ws.Range("M69").Value = -8165
$ws.Range("N69").Value = -15275
# Row 72
$ws.Range("H72").Value = 3761
$ws.Range("I72").Value = 3013
$ws.Range("J72").Value = 4509
$ws.Range("K72").Value = 27117
$ws.Range("L72").Value = 40581
$ws.Range("M72").Value = -22749
$ws.Range("N72").Value = -49317
# Row 88
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 3000
$ws.Range("N88").Value = -3812
# Row 91
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 3000
$ws.Range("N91").Value = -5808
# Row 98
$ws.Range("H98").Value = 1506.4
$ws.Range("I98").Value = 1650
$ws.Range("J98").Value = 1342.2858
$ws.Range("K98").Value = 1650
$ws.Range("L98").Value = 1342.2858
$ws.Range("M98").Value = -152
$ws.Range("N98").Value = -4338.2858
# Row 122
$ws.Range("H122").Value = 1506.4
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 1342.2858
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 4026.8574
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -8926.857400000001
# Row 135
$ws.Range("H135").Value = 558.087
$ws.Range("I135").Value = 324.33334
$ws.Range("J135").Value = 1399.6
$ws.Range("K135").Value = 2919.00006
$ws.Range("L135").Value = 12596.4
$ws.Range("M135").Value = -384.0000600000003
$ws.Range("N135").Value = -17666.4

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 7217.7144
$ws.Range("I74").Value = 9482.714
$ws.Range("J74").Value = 2687.7144
$ws.Range("K74").Value = 9482.714
$ws.Range("L74").Value = 2687.7144
$ws.Range("M74").Value = -8608.714
$ws.Range("N74").Value = -4435.7144
# Row 77
$ws.Range("H77").Value = 7217.7144
$ws.Range("I77").Value = 9482.714
$ws.Range("J77").Value = 2687.7144
$ws.Range("K77").Value = 47413.57
$ws.Range("L77").Value = 13438.572
$ws.Range("M77").Value = -43045.57
$ws.Range("N77").Value = -22174.572

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3518.9167
$ws.Range("I31").Value = 1920.9062
$ws.Range("J31").Value = 4502.3076
$ws.Range("K31").Value = 1920.9062
$ws.Range("L31").Value = 4502.3076
$ws.Range("M31").Value = -1625.9062
$ws.Range("N31").Value = -5092.3076
# Row 34
$ws.Range("H34").Value = 3518.9167
$ws.Range("I34").Value = 1920.9062
$ws.Range("J34").Value = 4502.3076
$ws.Range("K34").Value = 1920.9062
$ws.Range("L34").Value = 4502.3076
$ws.Range("M34").Value = -1718.9062
$ws.Range("N34").Value = -4906.3076

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -286
# Row 5
$ws.Range("H5").Value = 878217.1
$ws.Range("I5").Value = 505.1875
$ws.Range("J5").Value = 4389065
$ws.Range("K5").Value = 1515.5625
$ws.Range("L5").Value = 13167195
$ws.Range("M5").Value = -1403.5625
$ws.Range("N5").Value = -13167419
# Row 68
$ws.Range("H68").Value = 842932.9399999999
$ws.Range("I68").Value = 1231482.8
$ws.Range("J68").Value = 1075
$ws.Range("K68").Value = 3694448.4
$ws.Range("L68").Value = 3225
$ws.Range("M68").Value = -3693637.4
$ws.Range("N68").Value = -4847
# Row 70
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = 9000
$ws.Range("N70").Value = -9630
# Row 71
$ws.Range("H71").Value = 842932.9399999999
$ws.Range("I71").Value = 1231482.8
$ws.Range("J71").Value = 1075
$ws.Range("K71").Value = 11083345.2
$ws.Range("L71").Value = 9675
$ws.Range("M71").Value = -11079289.2
$ws.Range("N71").Value = -17787
# Row 73
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = 9000
$ws.Range("N73").Value = -11184
# Row 80
$ws.Range("H80").Value = 4449.5
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4561
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 13683
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -15555
# Row 83
$ws.Range("H83").Value = 4449.5
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4561
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 41049
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -50409
# Row 122
$ws.Range("H122").Value = 1213.95
$ws.Range("I122").Value = 532
$ws.Range("J122").Value = 3259.8
$ws.Range("K122").Value = 4788
$ws.Range("L122").Value = 29338.2
$ws.Range("M122").Value = -2338
$ws.Range("N122").Value = -34238.2
# Row 135
$ws.Range("H135").Value = 878217.1
$ws.Range("I135").Value = 505.1875
$ws.Range("J135").Value = 4389065
$ws.Range("K135").Value = 4546.6875
$ws.Range("L135").Value = 39501585
$ws.Range("M135").Value = -2011.6875
$ws.Range("N135").Value = -39506655

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 37830.8
$ws.Range("I40").Value = 53246.2
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 53246.2
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -53110.2
$ws.Range("N40").Value = -7272
# Row 100
$ws.Range("H100").Value = 999.6667
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 999.6667
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").Value = 999.6667
$ws.Range("N100").Value = -2081.6667
# Row 132
$ws.Range("H132").Value = 6276.1646
$ws.Range("I132").Value = 8503
$ws.Range("J132").Value = 4334.8203
$ws.Range("K132").Value = 25509
$ws.Range("L132").Value = 13004.4609
$ws.Range("M132").Value = -22979
$ws.Range("N132").Value = -18064.4609

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1957.6666
$ws.Range("I113").Value = 2774
$ws.Range("J113").Value = 325
$ws.Range("K113").Value = 8322
$ws.Range("L113").Value = 975
$ws.Range("M113").Value = -6152
$ws.Range("N113").Value = -5315
